$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '26.899.36'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '1.643.63'
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("E4").Value = '  -0.74%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '216.77'
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("E6").Value = '  +1.81%  '
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("E8").Value = '  +1.73%  '
$ws.Range("E9").Value = '  +0.51%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '19.83'
$ws.Range("E10").Value = '  +4.69%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.0846'
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '1.873.70'
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '1.626.47'
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '4.13'
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.528'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '66.27'
$ws.Range("E16").Value = '  +3.56%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '26.904.21'
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("E18").Value = '  +0.78%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '219.91'
$ws.Range("E19").Value = '  +4.22%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '4.38'
$ws.Range("E21").Value = '  +1.69%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '6.62'
$ws.Range("E22").Value = '  +7.12%  '
$ws.Range("E23").Value = '  +3.73%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '9.18'
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '145.79'
$ws.Range("E25").Value = '  -0.66%  '
$ws.Range("E26").Value = '  -0.76%  '
$ws.Range("E27").Value = '  +6.07%  '
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '15.82'
$ws.Range("E29").Value = '  +2.01%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '0.0505'
$ws.Range("E30").Value = '  +0.97%  '
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("E33").Value = '  +2.24%  '
$ws.Range("E34").Value = '  +2.79%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '1.246.47'
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '2.44'
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("E37").Value = '  +1.28%  '
$ws.Range("E38").Value = '  +2.90%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '0.834'
$ws.Range("E39").Value = '  +3.97%  '
$ws.Range("E40").Value = '  -0.67%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '0.806'
$ws.Range("E41").Value = '  +1.08%  '
$ws.Range("E42").Value = '  +2.25%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '1.785.05'
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '2.10'
$ws.Range("E44").Value = '  -2.39%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '60.73'
$ws.Range("E45").Value = '  +1.62%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '91.48'
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("E47").Value = '  +0.99%  '
$ws.Range("E48").Value = '  +11.94%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '0.0515'
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '0.0973'
$ws.Range("E50").Value = '  +1.98%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '7.57'
$ws.Range("E51").Value = '  +1.71%  '
